$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value2 = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMMF53DSW%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T160033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJHMEUCIAJRH%2Fww2VE%2BBaoUVyVI1aHA0jhEsAcpL6cWC4uYye0DAiEA%2B21Z9Pa%2BsqKydr2y5JmtxXQlUK%2BQPooAlc%2F6WJTM8%2BEq5gMIi%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBMIEDwwJZL73D%2BIrSq6Ax0A5xMrk%2F8ksq0FfXkZGmFn3fHJiWV2f67SOWX0isDfHOfHHx97md5iTgkOpzCvZVHEXmfKu1%2F3yj0gRMiAJcTgc%2BHHBXCsSSX0E1KMhc7YmSV%2FsPhqr0AxiPWlIvD%2BE8vmBVinxY%2BqaqHkBCXjqOKYQoiQqGoQxKLsQ08npr8zXzOPjgdxNRlFay8YxP37r4AgiyO2vcN4y%2FloVuHeLqeciuqx0AXbQrGYn9SnwfiMqYDKgiLGQtO%2FplDfcOwRVYLi3e9dl7jAVTGfBR5xvMZSwk69wxhhQGiT7Mc5G913edR2MXKLN%2F%2BZNvy1QfANc9Vq1kbv38zabm9G605Kig8C8hjrB0sulAI8J5rzpXB6vCuX2H6ASEYYVfLGt5OsnT31y4n8rNWLvgnY0eNh6a7zNwDJlbhcSBFa68UIRpaimj1rcH5zPsgvXPG%2F8%2BfcMQO1TizGZbGW1XpduaJDt8PLHCkSAwg645E%2B0%2BRxvXQkgsIStgqZqknUBYx12U8ozgn%2FmHucHWChVt8gWawKJhXL7fSsKCoifcWtZHWMcEqR3hbn4fmvtw0yrxecye02JsoMjzzI5lck8row36%2BkvwY6pQHyb5MVjyahweUSASjJ2TfolIrniaI2uE31x1PbIYwmBzx%2Bt8TkPuxxS9Jnou2g0u94cq0A%2FiUU7X8PYaqGfCg5%2F7ceJrlN%2Fi8p2BxW0oTIITgYZWlJUVR4HOqQv7zjp3KPV8BVHSjeDAuYKs7QZ2v12hQHUvRsIZRNs7xL6R7sKttQmPDamZRjV2Nt2nmefp23Y5oHOx%2B2cpQWcZgvw2xRI7F4aSM%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2000%3A00%3A33%20GMT&X-Amz-Signature=bb7b43974db680031dbe63a67a5fc868094e03701be498cd53c67bb43038bb20'
$ws.Range("K3").Value2 = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMMF53DSW%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T160033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJHMEUCIAJRH%2Fww2VE%2BBaoUVyVI1aHA0jhEsAcpL6cWC4uYye0DAiEA%2B21Z9Pa%2BsqKydr2y5JmtxXQlUK%2BQPooAlc%2F6WJTM8%2BEq5gMIi%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBMIEDwwJZL73D%2BIrSq6Ax0A5xMrk%2F8ksq0FfXkZGmFn3fHJiWV2f67SOWX0isDfHOfHHx97md5iTgkOpzCvZVHEXmfKu1%2F3yj0gRMiAJcTgc%2BHHBXCsSSX0E1KMhc7YmSV%2FsPhqr0AxiPWlIvD%2BE8vmBVinxY%2BqaqHkBCXjqOKYQoiQqGoQxKLsQ08npr8zXzOPjgdxNRlFay8YxP37r4AgiyO2vcN4y%2FloVuHeLqeciuqx0AXbQrGYn9SnwfiMqYDKgiLGQtO%2FplDfcOwRVYLi3e9dl7jAVTGfBR5xvMZSwk69wxhhQGiT7Mc5G913edR2MXKLN%2F%2BZNvy1QfANc9Vq1kbv38zabm9G605Kig8C8hjrB0sulAI8J5rzpXB6vCuX2H6ASEYYVfLGt5OsnT31y4n8rNWLvgnY0eNh6a7zNwDJlbhcSBFa68UIRpaimj1rcH5zPsgvXPG%2F8%2BfcMQO1TizGZbGW1XpduaJDt8PLHCkSAwg645E%2B0%2BRxvXQkgsIStgqZqknUBYx12U8ozgn%2FmHucHWChVt8gWawKJhXL7fSsKCoifcWtZHWMcEqR3hbn4fmvtw0yrxecye02JsoMjzzI5lck8row36%2BkvwY6pQHyb5MVjyahweUSASjJ2TfolIrniaI2uE31x1PbIYwmBzx%2Bt8TkPuxxS9Jnou2g0u94cq0A%2FiUU7X8PYaqGfCg5%2F7ceJrlN%2Fi8p2BxW0oTIITgYZWlJUVR4HOqQv7zjp3KPV8BVHSjeDAuYKs7QZ2v12hQHUvRsIZRNs7xL6R7sKttQmPDamZRjV2Nt2nmefp23Y5oHOx%2B2cpQWcZgvw2xRI7F4aSM%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2000%3A00%3A33%20GMT&X-Amz-Signature=565b625fbe71d034c397fbe7a10c1dbe3311800177d44668448b4106202f94bc'
$ws.Range("K4").Value2 = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMMF53DSW%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T160033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJHMEUCIAJRH%2Fww2VE%2BBaoUVyVI1aHA0jhEsAcpL6cWC4uYye0DAiEA%2B21Z9Pa%2BsqKydr2y5JmtxXQlUK%2BQPooAlc%2F6WJTM8%2BEq5gMIi%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBMIEDwwJZL73D%2BIrSq6Ax0A5xMrk%2F8ksq0FfXkZGmFn3fHJiWV2f67SOWX0isDfHOfHHx97md5iTgkOpzCvZVHEXmfKu1%2F3yj0gRMiAJcTgc%2BHHBXCsSSX0E1KMhc7YmSV%2FsPhqr0AxiPWlIvD%2BE8vmBVinxY%2BqaqHkBCXjqOKYQoiQqGoQxKLsQ08npr8zXzOPjgdxNRlFay8YxP37r4AgiyO2vcN4y%2FloVuHeLqeciuqx0AXbQrGYn9SnwfiMqYDKgiLGQtO%2FplDfcOwRVYLi3e9dl7jAVTGfBR5xvMZSwk69wxhhQGiT7Mc5G913edR2MXKLN%2F%2BZNvy1QfANc9Vq1kbv38zabm9G605Kig8C8hjrB0sulAI8J5rzpXB6vCuX2H6ASEYYVfLGt5OsnT31y4n8rNWLvgnY0eNh6a7zNwDJlbhcSBFa68UIRpaimj1rcH5zPsgvXPG%2F8%2BfcMQO1TizGZbGW1XpduaJDt8PLHCkSAwg645E%2B0%2BRxvXQkgsIStgqZqknUBYx12U8ozgn%2FmHucHWChVt8gWawKJhXL7fSsKCoifcWtZHWMcEqR3hbn4fmvtw0yrxecye02JsoMjzzI5lck8row36%2BkvwY6pQHyb5MVjyahweUSASjJ2TfolIrniaI2uE31x1PbIYwmBzx%2Bt8TkPuxxS9Jnou2g0u94cq0A%2FiUU7X8PYaqGfCg5%2F7ceJrlN%2Fi8p2BxW0oTIITgYZWlJUVR4HOqQv7zjp3KPV8BVHSjeDAuYKs7QZ2v12hQHUvRsIZRNs7xL6R7sKttQmPDamZRjV2Nt2nmefp23Y5oHOx%2B2cpQWcZgvw2xRI7F4aSM%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2000%3A00%3A33%20GMT&X-Amz-Signature=acc1863b1db5bdcf368016a7fac7477b522784c6c8a606eb3bc2fb73d3efb6e7'
$ws.Range("K5").Value2 = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMMF53DSW%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T160033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJHMEUCIAJRH%2Fww2VE%2BBaoUVyVI1aHA0jhEsAcpL6cWC4uYye0DAiEA%2B21Z9Pa%2BsqKydr2y5JmtxXQlUK%2BQPooAlc%2F6WJTM8%2BEq5gMIi%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBMIEDwwJZL73D%2BIrSq6Ax0A5xMrk%2F8ksq0FfXkZGmFn3fHJiWV2f67SOWX0isDfHOfHHx97md5iTgkOpzCvZVHEXmfKu1%2F3yj0gRMiAJcTgc%2BHHBXCsSSX0E1KMhc7YmSV%2FsPhqr0AxiPWlIvD%2BE8vmBVinxY%2BqaqHkBCXjqOKYQoiQqGoQxKLsQ08npr8zXzOPjgdxNRlFay8YxP37r4AgiyO2vcN4y%2FloVuHeLqeciuqx0AXbQrGYn9SnwfiMqYDKgiLGQtO%2FplDfcOwRVYLi3e9dl7jAVTGfBR5xvMZSwk69wxhhQGiT7Mc5G913edR2MXKLN%2F%2BZNvy1QfANc9Vq1kbv38zabm9G605Kig8C8hjrB0sulAI8J5rzpXB6vCuX2H6ASEYYVfLGt5OsnT31y4n8rNWLvgnY0eNh6a7zNwDJlbhcSBFa68UIRpaimj1rcH5zPsgvXPG%2F8%2BfcMQO1TizGZbGW1XpduaJDt8PLHCkSAwg645E%2B0%2BRxvXQkgsIStgqZqknUBYx12U8ozgn%2FmHucHWChVt8gWawKJhXL7fSsKCoifcWtZHWMcEqR3hbn4fmvtw0yrxecye02JsoMjzzI5lck8row36%2BkvwY6pQHyb5MVjyahweUSASjJ2TfolIrniaI2uE31x1PbIYwmBzx%2Bt8TkPuxxS9Jnou2g0u94cq0A%2FiUU7X8PYaqGfCg5%2F7ceJrlN%2Fi8p2BxW0oTIITgYZWlJUVR4HOqQv7zjp3KPV8BVHSjeDAuYKs7QZ2v12hQHUvRsIZRNs7xL6R7sKttQmPDamZRjV2Nt2nmefp23Y5oHOx%2B2cpQWcZgvw2xRI7F4aSM%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2000%3A00%3A33%20GMT&X-Amz-Signature=c28054bc020bc5f23809ce78a4b82fba8baace3749927e724bd4fb10b6038004'
$ws.Range("K6").Value2 = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMMF53DSW%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T160033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECIaCXVzLXdlc3QtMiJHMEUCIAJRH%2Fww2VE%2BBaoUVyVI1aHA0jhEsAcpL6cWC4uYye0DAiEA%2B21Z9Pa%2BsqKydr2y5JmtxXQlUK%2BQPooAlc%2F6WJTM8%2BEq5gMIi%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDBMIEDwwJZL73D%2BIrSq6Ax0A5xMrk%2F8ksq0FfXkZGmFn3fHJiWV2f67SOWX0isDfHOfHHx97md5iTgkOpzCvZVHEXmfKu1%2F3yj0gRMiAJcTgc%2BHHBXCsSSX0E1KMhc7YmSV%2FsPhqr0AxiPWlIvD%2BE8vmBVinxY%2BqaqHkBCXjqOKYQoiQqGoQxKLsQ08npr8zXzOPjgdxNRlFay8YxP37r4AgiyO2vcN4y%2FloVuHeLqeciuqx0AXbQrGYn9SnwfiMqYDKgiLGQtO%2FplDfcOwRVYLi3e9dl7jAVTGfBR5xvMZSwk69wxhhQGiT7Mc5G913edR2MXKLN%2F%2BZNvy1QfANc9Vq1kbv38zabm9G605Kig8C8hjrB0sulAI8J5rzpXB6vCuX2H6ASEYYVfLGt5OsnT31y4n8rNWLvgnY0eNh6a7zNwDJlbhcSBFa68UIRpaimj1rcH5zPsgvXPG%2F8%2BfcMQO1TizGZbGW1XpduaJDt8PLHCkSAwg645E%2B0%2BRxvXQkgsIStgqZqknUBYx12U8ozgn%2FmHucHWChVt8gWawKJhXL7fSsKCoifcWtZHWMcEqR3hbn4fmvtw0yrxecye02JsoMjzzI5lck8row36%2BkvwY6pQHyb5MVjyahweUSASjJ2TfolIrniaI2uE31x1PbIYwmBzx%2Bt8TkPuxxS9Jnou2g0u94cq0A%2FiUU7X8PYaqGfCg5%2F7ceJrlN%2Fi8p2BxW0oTIITgYZWlJUVR4HOqQv7zjp3KPV8BVHSjeDAuYKs7QZ2v12hQHUvRsIZRNs7xL6R7sKttQmPDamZRjV2Nt2nmefp23Y5oHOx%2B2cpQWcZgvw2xRI7F4aSM%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2000%3A00%3A33%20GMT&X-Amz-Signature=30ecae64f1e0b91a5af65108e6782a51ac21f1eaa3e6d557c4eeb2de3c1f87ec'
